$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2").Value = "Đàn tranh 123"
$ws.Range("K2").Select()
